$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Push the existing row 6 entry (40918 date row) down to row 7, keeping its
# values and formatting intact.
$ws.Range("A7").Value = 40918
$ws.Range("B7").Value2 = $ws.Range("B6").Value2
$ws.Range("C7").Value2 = $ws.Range("C6").Value2
$ws.Range("D7").Value2 = $ws.Range("D6").Value2
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = $null
$ws.Range("G7").Value = $null
$ws.Range("H7").Value = $null

$ws.Range("A6:H6").Copy()
$ws.Range("A7:H7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Overwrite row 6 with the new report entry
$ws.Range("A6").Value = "24/9/2012"
$ws.Range("B6").Value = "chỉnh sửa hoàn thiện SRS cá nhân"
$ws.Range("C6").Value = "1phaanf SRS"
$ws.Range("D6").Value = "Hoàn thành"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = $null
$ws.Range("H6").Value = $null

# Row 6 gets a taller custom height to fit the new wrapped text
$ws.Range("A6:H6").RowHeight = 40.5

$ws.Range("F7").Select()
